# Update MSME definitions: replace the literal "<br/>" markup in the
# SMMLV definition strings with an actual newline character.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$newTexts = @{
    "C21" = "<=500 SMMLV `n(Salario Minimo Mensual Legal Vigente - Legal monthly minimum wage in force)"
    "C22" = ">500 - <=5,000 SMMLV `n(Salario Minimo Mensual  Legal Vigente - Legal monthly minimum wage in force)"
    "C23" = ">5,000 - <=30,000 SMMLV `n(Salario Minimo Mensual Legal  Vigente - Legal monthly minimum wage in force)"
    "C24" = ">30,000 SMMLV `n(Salario Minimo Mensual Legal  Vigente - Legal monthly minimum wage in force)"
    "C42" = "<=500 SMMLV `n(Salario Minimo Mensual Legal Vigente - Legal monthly minimum wage in force)"
    "C43" = ">500 - <=5,000 SMMLV `n(Salario Minimo Mensual  Legal Vigente - Legal monthly minimum wage in force)"
    "C44" = ">5,000 - <=30,000 SMMLV `n(Salario Minimo Mensual Legal  Vigente - Legal monthly minimum wage in force)"
    "C45" = ">30,000 SMMLV `n(Salario Minimo Mensual Legal  Vigente - Legal monthly minimum wage in force)"
}

foreach ($addr in $newTexts.Keys) {
    $ws.Range($addr).Value = $newTexts[$addr]
}
